# Kfin support and improvements
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fund names (column A) for the new data rows 2-11
$funds = @(
    "HDFC Balanced Advantage Fund Direct Plan Growth",
    "ICICI Prudential Nifty 50 Index Direct Plan Growth",
    "Kotak Emerging Equity Fund Direct Growth",
    "Parag Parikh Flexi Cap Fund Direct Growth",
    "SBI Bluechip Direct Plan Growth",
    "SBI Small Cap Fund Direct Growth",
    "Mirae Asset Emerging Bluechip Fund - Direct Plan Growth",
    "Axis Mid Cap Fund - Direct Growth",
    "NIPPON INDIA SMALL CAP FUND - DIRECT GROWTH PLAN GROWTH OPTION",
    "Overall"
)

# Numeric data for columns B..H for each of the 10 rows (rows 2-11)
$data = @(
    @(2999.574, 3132.127, 44.298, 4.419, 2.205, 42.302, 19.578),
    @(3999.799, 4055.674, 58.501, 1.397, 1.044, 9.023999999999999, 6.688),
    @(2999.846, 3113.145, 34.667, 3.777, 0.338, 47.207, 3.617),
    @(3999.803, 4031.681, 40.5, 0.797, 1.147, 7.399, 10.786),
    @(3999.815, 3952.529, 24, -1.182, -0.998, -16.576, -14.177),
    @(1999.897, 2032.165, 19.5, 1.613, -0.97, 34.73, -16.726),
    @(1999.843, 2020.721, 28, 1.044, -2.202, 14.498, -25.196),
    @(999.944, 1005.537, 25, 0.5590000000000001, -1.841, 8.484999999999999, -23.762),
    @(999.8920000000001, 1021.12, 25, 2.123, -1.841, 35.896, -23.762),
    @(23998.412, 24364.7, 36.412, 1.526, 0.099, 16.327, 0.995)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $funds[$i]
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, 2 + $c).Value = $vals[$c]
    }
}
